$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns (D:E), shifting old D:K to F:M
$ws.Range("D5:E102").Insert(-4161)

# Step 2: Copy number formats/styles from column F (which now holds the old column D formatting)
# into the two newly inserted blank columns D:E so they match (date / number styles).
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Step 3: Populate the new D/E column values row by row
$ws.Range("D7").Value2 = 43453
$ws.Range("E7").Value2 = 43341
$ws.Range("D8").Value2 = 102900
$ws.Range("E8").Value2 = 83900
$ws.Range("D9").Value2 = 34200
$ws.Range("E9").Value2 = 27400
$ws.Range("D10").Value2 = 68700
$ws.Range("E10").Value2 = 56500
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 1200
$ws.Range("E14").Value2 = 2200
$ws.Range("D15").Value2 = 4900
$ws.Range("E15").Value2 = 4100
$ws.Range("D17").Value2 = 108600
$ws.Range("E17").Value2 = 84400
$ws.Range("D18").Value2 = -5700
$ws.Range("E18").Value2 = -500
$ws.Range("D20").Value2 = 0
$ws.Range("E20").Value2 = 0
$ws.Range("D21").Value2 = -700
$ws.Range("E21").Value2 = 3500
$ws.Range("D22").Value2 = 1700
$ws.Range("E22").Value2 = 1100
$ws.Range("D23").Value2 = -7400
$ws.Range("E23").Value2 = -1600
$ws.Range("D24").Value2 = 100
$ws.Range("E24").Value2 = 200
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = -7500
$ws.Range("E26").Value2 = -1900
$ws.Range("D27").Value2 = -7500
$ws.Range("E27").Value2 = -1900
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("E32").Value2 = 0
$ws.Range("D33").Value2 = -7500
$ws.Range("E33").Value2 = -1900
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = -7500
$ws.Range("E35").Value2 = -1900
$ws.Range("D38").Value2 = 43453
$ws.Range("E38").Value2 = 43341
$ws.Range("D41").Value2 = 8700
$ws.Range("E41").Value2 = 3700
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 8100
$ws.Range("E43").Value2 = 8800
$ws.Range("D44").Value2 = 4100
$ws.Range("E44").Value2 = 4000
$ws.Range("D45").Value2 = 12700
$ws.Range("E45").Value2 = 3200
$ws.Range("D46").Value2 = 33500
$ws.Range("E46").Value2 = 19800
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 133600
$ws.Range("E48").Value2 = 138300
$ws.Range("D49").Value2 = 18200
$ws.Range("E49").Value2 = 18700
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 23500
$ws.Range("E52").Value2 = 23200
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 208900
$ws.Range("E54").Value2 = 200000
$ws.Range("D57").Value2 = 8300
$ws.Range("E57").Value2 = 10500
$ws.Range("D58").Value2 = 10000
$ws.Range("E58").Value2 = 39300
$ws.Range("D59").Value2 = 30800
$ws.Range("E59").Value2 = 31800
$ws.Range("D60").Value2 = 49000
$ws.Range("E60").Value2 = 81600
$ws.Range("D61").Value2 = 46100
$ws.Range("E61").Value2 = 0
$ws.Range("D62").Value2 = 5700
$ws.Range("E62").Value2 = 5800
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 100800
$ws.Range("E66").Value2 = 87400
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 68900
$ws.Range("E72").Value2 = 73900
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 108100
$ws.Range("E76").Value2 = 112600
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43453
$ws.Range("E80").Value2 = 43341
$ws.Range("D81").Value2 = -7500
$ws.Range("E81").Value2 = -1900
$ws.Range("D83").Value2 = 4900
$ws.Range("E83").Value2 = 4100
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 1300
$ws.Range("E89").Value2 = -3600
$ws.Range("D91").Value2 = -1100
$ws.Range("E91").Value2 = -1500
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -900
$ws.Range("E94").Value2 = 10600
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = 15700
$ws.Range("E100").Value2 = -4800
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("D102").Value2 = 16000
$ws.Range("E102").Value2 = 2200

# Step 4: Apply value corrections to a handful of existing cells (now shifted into columns G/H)
$ws.Range("H8").Value2 = 113500
$ws.Range("H10").Value2 = 76900
$ws.Range("H18").Value2 = -5000
$ws.Range("H20").Value2 = 100
$ws.Range("H21").Value2 = 500
$ws.Range("H23").Value2 = -5500
$ws.Range("G24").Value2 = 200
$ws.Range("G26").Value2 = -13800
$ws.Range("H26").Value2 = -5500
$ws.Range("G27").Value2 = -13800
$ws.Range("H27").Value2 = -5500
$ws.Range("G29").Value2 = -3300
$ws.Range("H32").Value2 = -100
$ws.Range("H33").Value2 = -5500
$ws.Range("H35").Value2 = -5500
$ws.Range("H81").Value2 = -5500
